# Auto-generated edit script: updates title date and all 100 table cell expressions
$d = $word.ActiveDocument

# Update title paragraph (date line)
$d.Paragraphs.Item(1).Range.Text = "2026-01-17 Saturday"

# Update each table cell with its new arithmetic expression, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "0+17=",
    "62-48=",
    "81-34=",
    "1+62=",
    "88+7=",
    "44-10=",
    "66+19=",
    "95-69=",
    "86-75=",
    "32+58=",
    "53+6=",
    "31-12=",
    "27-18=",
    "1+69=",
    "94-60=",
    "39+14=",
    "16+46=",
    "86-55=",
    "5+58=",
    "10+63=",
    "94-47=",
    "43-25=",
    "23+41=",
    "37+46=",
    "55+36=",
    "69-28=",
    "96-40=",
    "28+60=",
    "74-63=",
    "17+64=",
    "39+9=",
    "88-57=",
    "32+20=",
    "5+91=",
    "6+19=",
    "93-12=",
    "45+18=",
    "89-56=",
    "62-21=",
    "83-31=",
    "62+20=",
    "87-13=",
    "86-85=",
    "86-77=",
    "51+14=",
    "56+6=",
    "94+3=",
    "22+18=",
    "5+71=",
    "32+48=",
    "82-55=",
    "83-8=",
    "86-25=",
    "86-13=",
    "36+29=",
    "0+5=",
    "5+37=",
    "49-20=",
    "80+1=",
    "66+18=",
    "48+12=",
    "38+45=",
    "22+20=",
    "74-57=",
    "76-16=",
    "98-32=",
    "69-67=",
    "21+69=",
    "63-52=",
    "68-24=",
    "1+98=",
    "15+57=",
    "0+2=",
    "55-9=",
    "30+51=",
    "53+11=",
    "52+24=",
    "2+10=",
    "12+6=",
    "79-38=",
    "47-47=",
    "61-19=",
    "20+70=",
    "95-6=",
    "31+44=",
    "15+49=",
    "3+57=",
    "87-15=",
    "86-6=",
    "52-33=",
    "61+29=",
    "34-20=",
    "91-23=",
    "58+12=",
    "73-1=",
    "77-60=",
    "73-15=",
    "56-27=",
    "23-16=",
    "74+17="
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx++
    }
}
